$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets('Citywide Totals')
$ws.Range('D2').Value = 67
$ws.Range('E2').Value = 51
$ws.Range('I2').Value = 87
$ws.Range('K2').Value = 102
$ws.Range('J3').Value = 151
$ws.Range('B9').Value = 275
$ws.Range('C9').Value = 345
$ws.Range('E9').Value = 315
$ws.Range('I9').Value = 383
$ws.Range('K9').Value = 372
$ws.Range('B10').Value = 919
$ws.Range('C10').Value = 1107
$ws.Range('E10').Value = 1570
$ws.Range('F10').Value = 1605
$ws.Range('G10').Value = 769
$ws.Range('H10').Value = 394
$ws.Range('J10').Value = 513
$ws.Range('K10').Value = 523
$ws.Range('B11').Value = 1299
$ws.Range('C11').Value = 1574
$ws.Range('D11').Value = 1745
$ws.Range('E11').Value = 2046
$ws.Range('F11').Value = 2160
$ws.Range('G11').Value = 1286
$ws.Range('H11').Value = 903
$ws.Range('I11').Value = 1269
$ws.Range('J11').Value = 1066
$ws.Range('K11').Value = 1180

$ws = $wb.Worksheets('Chinatown')
$ws.Range('K7').Value = 16
$ws.Range('K9').Value = 26

$ws = $wb.Worksheets('Garfield Park')
$ws.Range('I2').Value = 4
$ws.Range('C7').Value = 23
$ws.Range('E8').Value = 58
$ws.Range('C9').Value = 66
$ws.Range('E9').Value = 96
$ws.Range('I9').Value = 74

$ws = $wb.Worksheets('Grand Crossing')
$ws.Range('B7').Value = 15
$ws.Range('C8').Value = 43
$ws.Range('B9').Value = 47
$ws.Range('C9').Value = 72

$ws = $wb.Worksheets('Loop')
$ws.Range('C8').Value = 27
$ws.Range('E8').Value = 49
$ws.Range('B9').Value = 126
$ws.Range('C9').Value = 213
$ws.Range('E9').Value = 453
$ws.Range('F9').Value = 429
$ws.Range('G9').Value = 142
$ws.Range('B10').Value = 162
$ws.Range('C10').Value = 255
$ws.Range('E10').Value = 513
$ws.Range('F10').Value = 484
$ws.Range('G10').Value = 205

$ws = $wb.Worksheets('Old Town')
$ws.Range('D2').Value = 2
$ws.Range('I5').Value = 14
$ws.Range('B6').Value = 17
$ws.Range('B7').Value = 19
$ws.Range('D7').Value = 37
$ws.Range('I7').Value = 26

$ws = $wb.Worksheets('Little Italy, UIC')
$ws.Range('E5').Value = 12
$ws.Range('B6').Value = 14
$ws.Range('K6').Value = 9
$ws.Range('B7').Value = 22
$ws.Range('E7').Value = 33
$ws.Range('K7').Value = 33

$ws = $wb.Worksheets('North Lawndale')
$ws.Range('J3').Value = 2
$ws.Range('E7').Value = 24
$ws.Range('E8').Value = 32
$ws.Range('J8').Value = 9

$ws = $wb.Worksheets('By Neighborhood')
$ws.Range('C2').Value = 9
$ws.Range('I2').Value = 6
$ws.Range('I7').Value = 16
$ws.Range('K21').Value = 26
$ws.Range('E27').Value = 21
$ws.Range('B28').Value = 74
$ws.Range('H28').Value = 57
$ws.Range('I28').Value = 65
$ws.Range('K28').Value = 65
$ws.Range('E29').Value = 18
$ws.Range('C32').Value = 66
$ws.Range('E32').Value = 96
$ws.Range('I32').Value = 74
$ws.Range('B36').Value = 47
$ws.Range('C36').Value = 72
$ws.Range('K41').Value = 13
$ws.Range('B43').Value = 8
$ws.Range('D47').Value = 40
$ws.Range('E47').Value = 51
$ws.Range('K47').Value = 29
$ws.Range('B50').Value = 22
$ws.Range('E50').Value = 33
$ws.Range('K50').Value = 33
$ws.Range('B53').Value = 162
$ws.Range('C53').Value = 255
$ws.Range('E53').Value = 513
$ws.Range('F53').Value = 484
$ws.Range('G53').Value = 205
$ws.Range('K54').Value = 7
$ws.Range('B56').Value = 6
$ws.Range('E65').Value = 32
$ws.Range('J65').Value = 9
$ws.Range('F68').Value = 32
$ws.Range('B70').Value = 19
$ws.Range('D70').Value = 37
$ws.Range('I70').Value = 26
$ws.Range('C74').Value = 29
$ws.Range('F74').Value = 75
$ws.Range('J74').Value = 28
$ws.Range('G76').Value = 34
$ws.Range('K77').Value = 48
$ws.Range('E78').Value = 38
$ws.Range('B81').Value = 21
$ws.Range('E87').Value = 26
$ws.Range('B92').Value = 18
$ws.Range('K95').Value = 21
$ws.Range('B99').Value = 1299
$ws.Range('C99').Value = 1574
$ws.Range('D99').Value = 1745
$ws.Range('E99').Value = 2046
$ws.Range('F99').Value = 2160
$ws.Range('G99').Value = 1286
$ws.Range('H99').Value = 903
$ws.Range('I99').Value = 1269
$ws.Range('J99').Value = 1066
$ws.Range('K99').Value = 1180

$ws = $wb.Worksheets('Humboldt Park')
$ws.Range('K6').Value = 2
$ws.Range('K7').Value = 13

$ws = $wb.Worksheets('South Chicago')
$ws.Range('B5').Value = 9
$ws.Range('B6').Value = 21

$ws = $wb.Worksheets('Uptown')
$ws.Range('E8').Value = 18
$ws.Range('E9').Value = 26

$ws = $wb.Worksheets('Rush & Division')
$ws.Range('E2').Value = 1
$ws.Range('E6').Value = 38

$ws = $wb.Worksheets('Englewood')
$ws.Range('I2').Value = 7
$ws.Range('K7').Value = 26
$ws.Range('B8').Value = 46
$ws.Range('H8').Value = 20
$ws.Range('B9').Value = 74
$ws.Range('H9').Value = 57
$ws.Range('I9').Value = 65
$ws.Range('K9').Value = 65

$ws = $wb.Worksheets('Lake View')
$ws.Range('D2').Value = 1
$ws.Range('K2').Value = 6
$ws.Range('E7').Value = 40
$ws.Range('D8').Value = 40
$ws.Range('E8').Value = 51
$ws.Range('K8').Value = 29

$ws = $wb.Worksheets('Fuller Park')
$ws.Range('E8').Value = 13
$ws.Range('E9').Value = 18

$ws = $wb.Worksheets('Rogers Park')
$ws.Range('G7').Value = 23
$ws.Range('G8').Value = 34

$ws = $wb.Worksheets('River North')
$ws.Range('J3').Value = 3
$ws.Range('C6').Value = 24
$ws.Range('F6').Value = 65
$ws.Range('J6').Value = 9
$ws.Range('C7').Value = 29
$ws.Range('F7').Value = 75
$ws.Range('J7').Value = 28

$ws = $wb.Worksheets('West Loop')
$ws.Range('B8').Value = 12
$ws.Range('B9').Value = 18

$ws = $wb.Worksheets('Edgewater')
$ws.Range('E6').Value = 15
$ws.Range('E7').Value = 21

$ws = $wb.Worksheets('Albany Park')
$ws.Range('I2').Value = 1
$ws.Range('C6').Value = 8
$ws.Range('C7').Value = 9
$ws.Range('I7').Value = 6

$ws = $wb.Worksheets('Lower West Side')
$ws.Range('K5').Value = 4
$ws.Range('K6').Value = 7

$ws = $wb.Worksheets('Roseland')
$ws.Range('K9').Value = 27
$ws.Range('K10').Value = 48

$ws = $wb.Worksheets('West Town')
$ws.Range('K6').Value = 13
$ws.Range('K7').Value = 21

$ws = $wb.Worksheets('Auburn Gresham')
$ws.Range('I5').Value = 4
$ws.Range('I7').Value = 16

$ws = $wb.Worksheets('O''Hare')
$ws.Range('F6').Value = 30
$ws.Range('F7').Value = 32

$ws = $wb.Worksheets('Avondale')
$ws.Range('B6').Value = 7

$ws = $wb.Worksheets('Irving Park')
$ws.Range('B7').Value = 8

$ws = $wb.Worksheets('Mckinley Park')
$ws.Range('B5').Value = 5
$ws.Range('B6').Value = 6
